$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Header row (row 3) -----
$ws.Range("A3").Value = "No"
$ws.Range("B3").Value = "API"
$ws.Range("C3").Value = "Description"
$ws.Range("D3").Value = "Input"
$ws.Range("E3").Value = "Result"
$ws.Range("F3").Value = "possible values"
$ws.Range("G3").Value = "detail"

# ----- Row 4: sign-up -----
$ws.Range("A4").Value = 1
$ws.Range("C4").Value = "Sign Up new user"
$ws.Range("D4").Value = @"

"name": "Tiny",
"email": "tiny9611@outlook.com",
"phone": "+8613522171058"

"@
$ws.Range("E4").Value = @"

success: true/false
data {
message:“invalid_email”
}

"@
$ws.Range("F4").Value = @"
-invalid_email
-invalid_name
-invalid_phone
"@
$ws.Range("G4").Value = "send OTP codes through SMS or email"

# ----- Row 5: sign-in -----
$ws.Range("A5").Value = 2
$ws.Range("C5").Value = "Login "
$ws.Range("D5").Value = """phone"":""+8613522171058"""
$ws.Range("E5").Value = @"

success: true/false
data {
message:“invalid_email”
}

"@
$ws.Range("F5").Value = @"
-invalid_phone
-not_exising_user
"@
$ws.Range("G5").Value = "send OTP codes through SMS or email"

# ----- Row 6: check-verify-code -----
$ws.Range("A6").Value = 3
$ws.Range("C6").Value = "verify otp code"
$ws.Range("D6").Value = @"
"phone":"+8613522171058"
"otp": "4444"
"@
$ws.Range("E6").Value = @"

success: true/false
data {
message:“invalid_otp”
token:"123456"
}

"@
$ws.Range("F6").Value = @"
-invalid_phone
-invalid_otp
"@

# ----- Row 7: get-user-profile -----
$ws.Range("A7").Value = 4
$ws.Range("D7").Value = @"

"token":"123456

"@
$ws.Range("E7").Value = @"

success: true/false
data {
message:"invalid_otp"
profile: {
}
}

"@

# ----- Row 8: add-email -----
$ws.Range("A8").Value = 5
# D8 needs a leading apostrophe (quote-prefixed text)
$ws.Range("D8").Formula = "'-----"

# ----- Rows 9-15: sequence numbers only (B/C already had values) -----
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 7
$ws.Range("A11").Value = 8
$ws.Range("A12").Value = 9
$ws.Range("A13").Value = 10
$ws.Range("A14").Value = 11
$ws.Range("A15").Value = 12

# ----- Formatting -----
$ws.Range("C3:C6,C8,C9,C11:C13,D4:D7,E4:E7,G4:G5").WrapText = $true
$ws.Range("C18:C19").WrapText = $true
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F4:F6").NumberFormat = "@"
$ws.Range("F4:F6").WrapText = $true
$ws.Range("D8").WrapText = $true

# ----- Column widths -----
$ws.Columns.Item(4).ColumnWidth = 33.57
$ws.Columns.Item(5).ColumnWidth = 24.57
$ws.Columns.Item(6).ColumnWidth = 36.86

# ----- Row heights (wrapped multi-line cells) -----
$ws.Rows.Item(4).RowHeight = 105
$ws.Rows.Item(5).RowHeight = 105
$ws.Rows.Item(6).RowHeight = 120
$ws.Rows.Item(7).RowHeight = 150

# ----- Page setup -----
$ws.PageSetup.Orientation = 1

# ----- Selection / view -----
$ws.Range("E20").Select()

